# Restores the "Rules" sheet's R30 "From" threshold (cell C10) to 1.
# (revision restore / resave -- the only substantive data change vs. the
# prior revision is the C10 numeric value, 18 -> 1.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
